# Updates the cryptocurrency price/volume snapshot on Sheet1 (columns B-E,
# rows 2-51) to match the latest scrape, including a couple of rows whose
# coin order was swapped (WrappedEther/ShibaInu and Aptos/BabyDogeCoin).
#
# Column D ("Price") holds text-look-alike numbers (e.g. "240.19",
# "29.390.26") that must stay stored as text, exactly as they were before
# the edit. Assigning a plain string to .Value lets Excel auto-convert
# number-looking text into a real numeric cell, so for those cells we
# briefly force a text number format, set the value, then restore the
# cell's style to "Normal" so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = "29.390.26"
$c.Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  +0.01%  "
$ws.Cells.Item(3, 5).Value = "  -0.01%  "
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "240.19"
$c.Style = "Normal"
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "0.6298"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -0.86%  "
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = "0.9999"
$c.Style = "Normal"
$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = "0.07643"
$c.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +0.89%  "
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = "0.2933"
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  -1.03%  "
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "24.58"
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -1.05%  "
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "0.07737"
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +0.04%  "
$ws.Cells.Item(12, 2).Value = "ShibaInu"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "0.00001125"
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +12.80%  "
$ws.Cells.Item(13, 2).Value = "WrappedEther"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "1.859.81"
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -6.31%  "
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "5.005"
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -0.12%  "
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "0.6792"
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  -0.61%  "
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "83.71"
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +0.74%  "
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = "2.110.08"
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -6.81%  "
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = "6.185"
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +0.75%  "
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "29.407.57"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -0.03%  "
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = "228.78"
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -1.10%  "
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "12.47"
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +0.18%  "
$ws.Cells.Item(22, 5).Value = "  +0.05%  "
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "7.505"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -0.68%  "
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "157.36"
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +0.49%  "
$ws.Cells.Item(26, 5).Value = "  +0.00%  "
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "8.347"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -0.50%  "
$ws.Cells.Item(28, 5).Value = "  -0.37%  "
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "1.467"
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  -0.40%  "
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = "1.302"
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +3.76%  "
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "0.05598"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -1.78%  "
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "4.117"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -0.39%  "
$ws.Cells.Item(33, 5).Value = "  +0.21%  "
$ws.Cells.Item(34, 5).Value = "  +0.30%  "
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "1.157"
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +0.03%  "
$ws.Cells.Item(36, 5).Value = "  -0.71%  "
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = "2.587"
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -0.46%  "
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "1.242.05"
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -0.35%  "
$ws.Cells.Item(39, 5).Value = "  +0.03%  "
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "2.779"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -0.79%  "
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "6.421"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +5.38%  "
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "0.9026"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +0.11%  "
$ws.Cells.Item(43, 5).Value = "  -0.01%  "
$ws.Cells.Item(45, 5).Value = "  -0.28%  "
$ws.Cells.Item(46, 2).Value = "Aptos"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "7.159"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +1.25%  "
$ws.Cells.Item(47, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "0.00000000119"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +0.45%  "
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "0.4017"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -0.27%  "
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "9.037"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -1.15%  "
$ws.Cells.Item(50, 5).Value = "  -1.47%  "
$ws.Cells.Item(51, 5).Value = "  -0.43%  "

Write-Host "Applied all changes"